# Update the "Generate Report for Handback" timestamps.
# These cells hold plain text timestamps (not real Excel date serials),
# so we assign string values to keep them stored as shared-string text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for the first row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-07 01:26:16"

# --- zh-cn sheet: Correspond Handoff / Handback datetimes for the first row ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-07 01:26:11"
$wsZhCn.Range("K2").Value = "2016-09-07 01:26:37"

# --- de-de sheet: Correspond Handback DateTime for the first row ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-09-07 01:26:45"
